$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# In the "03.2-Stack-and-Queue-Advanced-Exercises" document there is a
# heading paragraph (style "3") that reads "00Насоки" because of a
# stray leftover run containing just "00" (with lang=bg-BG) sitting in
# front of the real "Насоки" run. Remove that stray "00" run so the
# heading just reads "Насоки".
$rng = $d.Content
$null = $rng.Find.Execute("00Насоки", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.SetRange($rng.Start, $rng.Start + 2)
    $rng.Delete()
}

# --- Edit 2 -----------------------------------------------------------
# Remove the trailing empty paragraph (its pPr only carries an
# rPr/lang="bg-BG") that sits right before the last section's sectPr,
# collapsing the document so the preceding empty paragraph is
# immediately followed by the section break.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$last = $d.Paragraphs.Last
$killRange = $d.Range($secondLast.Range.End - 1, $last.Range.End)
$killRange.Delete()
